$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JQUERY")
$ws.Activate()

# --- New "Test Data" values entered by the author, in the order they were typed ---
# (this order controls the order new entries are appended to the shared-strings table)

$ws.Range("F89").Value = "Faster"
$ws.Range("F91").Value = "ui.jQuery.js"
$ws.Range("F95").Value = "Dr."

$ws.Range("F14").Value = "XOffset = 100            Yoffset = 100"
$ws.Range("F14").WrapText = $true

$ws.Range("F24").Value = "XOffset = 200            Yoffset = 50"
$ws.Range("F24").WrapText = $true

# leading apostrophe forces text / quote-prefix (leading spaces preserved)
$ws.Range("F47").Value = "'       J"

$ws.Range("F59").Value = "    London"
$ws.Range("F61").Value = "    5 Star"

$ws.Range("F66").Value = "Automatic, Insurance, 1"
$ws.Range("F66").WrapText = $true

# "May/28/2021" would otherwise be auto-parsed as a date serial; force it to
# stay text (matching the source file) by typing it with a leading apostrophe
# into a scratch cell (outside the sheet's used range, pre-formatted to match
# F71's own style), then copying only the resulting value onto F71.
$scratch = $ws.Range("K71")
$ws.Range("F59").Copy()
$scratch.PasteSpecial(-4122)  # xlPasteFormats
$scratch.Value = "'May/28/2021"
$scratch.Copy()
$ws.Range("F71").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

$ws.Range("F78").Value = "Xoffset = 400, Yoffset = 50"
$ws.Range("F78").WrapText = $true

$ws.Range("F84").Value = "Music, Pop"

# --- numeric test data entries ---
$ws.Range("F93").Value = 3
$ws.Range("F100").Value = 400
$ws.Range("F104").Value = 12
$ws.Range("F113").Value = 45

# --- bold a handful of "Test Steps" labels ---
$ws.Range("E72").Font.Bold = $true
$ws.Range("E90").Font.Bold = $true
$ws.Range("E92").Font.Bold = $true
$ws.Range("E94").Font.Bold = $true

# --- restore view / selection state ---
$excel.ActiveWindow.ScrollRow = 48
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F84").Select()
